$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Using NumberFormat "@" before assignment and
# ClearFormats() after ensures numeric-looking strings (e.g. "1.000",
# "28.940.81") are stored as text (matching the source inline-string cells)
# instead of being auto-coerced to numbers, while keeping the original
# (default) cell style untouched.
$updates = @{
    'D2' = '28.940.81'
    'E2' = '  -1.52%  '
    'D3' = '1.833.19'
    'E3' = '  -1.90%  '
    'D4' = '1.000'
    'E4' = '  -0.11%  '
    'D5' = '245.41'
    'E5' = '  +0.66%  '
    'D6' = '0.6891'
    'E6' = '  -2.20%  '
    'D8' = '0.07691'
    'E8' = '  -2.91%  '
    'D9' = '0.3052'
    'E9' = '  -2.57%  '
    'D10' = '23.50'
    'E10' = '  -3.90%  '
    'D11' = '0.07820'
    'E11' = '  -0.41%  '
    'D12' = '1.834.19'
    'E12' = '  -1.66%  '
    'E13' = '  -1.95%  '
    'D14' = '90.41'
    'E14' = '  -3.63%  '
    'D15' = '0.6802'
    'E15' = '  -2.97%  '
    'D16' = '6.439'
    'E16' = '  -1.20%  '
    'D17' = '0.000008322'
    'E17' = '  -0.84%  '
    'D18' = '28.985.84'
    'E18' = '  -1.32%  '
    'D19' = '243.14'
    'E19' = '  -4.32%  '
    'D20' = '2.084.79'
    'E20' = '  -1.34%  '
    'E21' = '  -2.88%  '
    'D22' = '0.9997'
    'E22' = '  -0.01%  '
    'D23' = '7.477'
    'D24' = '1.001'
    'E24' = '  -0.07%  '
    'D25' = '163.28'
    'E25' = '  +1.38%  '
    'D26' = '0.1469'
    'E26' = '  -5.65%  '
    'D27' = '8.803'
    'E27' = '  -2.27%  '
    'D28' = '18.21'
    'E28' = '  -3.26%  '
    'E29' = '  +3.25%  '
    'D30' = '4.211'
    'E30' = '  -2.49%  '
    'D31' = '4.157'
    'E31' = '  -2.29%  '
    'D32' = '1.174'
    'E32' = '  -3.22%  '
    'D33' = '0.05113'
    'E33' = '  -3.16%  '
    'D34' = '0.7756'
    'E34' = '  +3.73%  '
    'D35' = '1.841'
    'E35' = '  -2.98%  '
    'D36' = '1.143'
    'E36' = '  -2.60%  '
    'D37' = '2.688'
    'E37' = '  -0.79%  '
    'D38' = '0.01848'
    'E38' = '  -1.59%  '
    'D39' = '1.242.45'
    'E39' = '  -2.88%  '
    'E40' = '  -2.46%  '
    'D41' = '0.9432'
    'E41' = '  +5.72%  '
    'D42' = '108.22'
    'E42' = '  -0.45%  '
    'D43' = '0.9995'
    'E43' = '  -0.14%  '
    'D44' = '5.695'
    'E44' = '  -5.06%  '
    'D45' = '9.585'
    'E45' = '  -0.16%  '
    'E46' = '  -5.17%  '
    'D47' = '1.984.75'
    'E47' = '  -1.58%  '
    'E48' = '  -0.10%  '
    'D49' = '64.54'
    'E49' = '  -9.03%  '
    'D50' = '1.748'
    'E50' = '  -2.78%  '
    'D51' = '0.4198'
    'E51' = '  -2.39%  '
}

foreach ($cell in $updates.Keys) {
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $updates[$cell]
    $range.ClearFormats()
}
